$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.741.14'
$ws.Range('E2').Value = '  -1.22%  '
$ws.Range('D3').Value = '2.028.20'
$ws.Range('E3').Value = '  -1.70%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.50'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.47%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.384'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0812'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('E12').Value = '  -1.04%  '
$ws.Range('D13').Value = '2.328.94'
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.00'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.756'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('D17').Value = '2.035.81'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').Value = '37.698.14'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.03'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.63'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').Value = '0.0₃0821'
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.81%  '
$ws.Range('E25').Value = '  -2.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.17'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('E28').Value = '  -3.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.83'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  -4.92%  '
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('E32').Value = '  -3.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.08'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.15%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.49'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.91%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0600'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.75%  '
$ws.Range('E37').Value = '  -3.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.23'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.36%  '
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').Value = '1.530.78'
$ws.Range('E40').Value = '  +3.28%  '
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '96.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('E45').Value = '  -3.34%  '
$ws.Range('E46').Value = '  -1.95%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.91'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.22%  '
$ws.Range('E48').Value = '  -1.74%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.14'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('D51').Value = '2.217.92'
$ws.Range('E51').Value = '  -1.65%  '
